$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$col = $ws.Columns.Item(1)

# --- Remove obsolete account rows -------------------------------------------------
# Each of these accounts had its row dropped entirely in the new export.
$toRemove = @("008004799", "004515548", "005141215", "004472538")
foreach ($acct in $toRemove) {
    $found = $col.Find($acct)
    if ($found -ne $null) {
        $ws.Rows.Item($found.Row).Delete()
    }
}

# --- Update the Saldo (balance) for accounts that kept their row -----------------
$updates = @{
    "004376853" = 42903.36
    "004499920" = 18730.93
    "004467884" = 10100
    "004479287" = 5497.82
}
foreach ($acct in $updates.Keys) {
    $found = $col.Find($acct)
    if ($found -ne $null) {
        $ws.Cells.Item($found.Row, 3).Value = $updates[$acct]
    }
}

# --- Insert the new account row right after 004467884 (ANA) ----------------------
$anchor = $col.Find("004467884")
$newRow = $anchor.Row + 1
$ws.Rows.Item($newRow).Insert()
$ws.Cells.Item($newRow, 1).Value = "'005198093"
$ws.Cells.Item($newRow, 2).Value = "ANA"
$ws.Cells.Item($newRow, 3).Value = 10000
